# Scheduled runner update: refresh market-board price snapshots and
# recomputed leve profit figures across the Anima Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 450
$ws.Range("I18").Value = 450
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 450
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -166

# Row 20
$ws.Range("H20").Value = 900
$ws.Range("I20").Value = 900
$ws.Range("K20").Value = 900
$ws.Range("M20").Value = -670

# Row 35
$ws.Range("H35").Value = 900
$ws.Range("I35").Value = 900
$ws.Range("K35").Value = 900
$ws.Range("M35").Value = -521

# Row 43
$ws.Range("H43").Value = 62511252
$ws.Range("I43").Value = 1475.25
$ws.Range("J43").Value = 125021020
$ws.Range("K43").Value = 1475.25
$ws.Range("L43").Value = 125021020
$ws.Range("M43").Value = -1406.25
$ws.Range("N43").Value = -125021158

# Row 88
$ws.Range("H88").Value = 101166.664
$ws.Range("J88").Value = 151000
$ws.Range("L88").Value = 151000
$ws.Range("N88").Value = -151812

# Row 91
$ws.Range("H91").Value = 101166.664
$ws.Range("J91").Value = 151000
$ws.Range("L91").Value = 151000
$ws.Range("N91").Value = -153808

# Row 132
$ws.Range("H132").Value = 1765.8549
$ws.Range("I132").Value = 1499.9122
$ws.Range("J132").Value = 4797.6
$ws.Range("K132").Value = 4499.7366
$ws.Range("L132").Value = 14392.8
$ws.Range("M132").Value = -1969.7366
$ws.Range("N132").Value = -19452.8

# Row 135
$ws.Range("H135").Value = 1236.075
$ws.Range("I135").Value = 488.59375
$ws.Range("J135").Value = 4226
$ws.Range("K135").Value = 4397.34375
$ws.Range("L135").Value = 38034
$ws.Range("M135").Value = -1862.34375
$ws.Range("N135").Value = -43104

# Row 138
$ws.Range("H138").Value = 2889.15
$ws.Range("I138").Value = 1554.0605
$ws.Range("J138").Value = 3826.5532
$ws.Range("K138").Value = 4662.181500000001
$ws.Range("L138").Value = 11479.6596
$ws.Range("M138").Value = 477.8184999999994
$ws.Range("N138").Value = -21759.6596

# Row 141
$ws.Range("H141").Value = 4874.0586
$ws.Range("I141").Value = 1982.6154
$ws.Range("J141").Value = 14271.25
$ws.Range("K141").Value = 5947.8462
$ws.Range("L141").Value = 42813.75
$ws.Range("M141").Value = -767.8462
$ws.Range("N141").Value = -53173.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 586746.7
$ws.Range("I32").Value = 650277.3
$ws.Range("J32").Value = 24046.857
$ws.Range("K32").Value = 650277.3
$ws.Range("L32").Value = 24046.857
$ws.Range("M32").Value = -649990.3
$ws.Range("N32").Value = -24620.857

# Row 61
$ws.Range("H61").Value = 2621.5818
$ws.Range("I61").Value = 2106.9412
$ws.Range("J61").Value = 3454.8096
$ws.Range("K61").Value = 2106.9412
$ws.Range("L61").Value = 3454.8096
$ws.Range("M61").Value = -1894.9412
$ws.Range("N61").Value = -3878.8096

# Row 74
$ws.Range("H74").Value = 2762.6155
$ws.Range("I74").Value = 2576.1667
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2576.1667
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1702.1667
$ws.Range("N74").Value = -6748

# Row 77
$ws.Range("H77").Value = 2762.6155
$ws.Range("I77").Value = 2576.1667
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 12880.8335
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -8512.833500000001
$ws.Range("N77").Value = -33736

# Row 132
$ws.Range("H132").Value = 2729.4666
$ws.Range("I132").Value = 1812.6046
$ws.Range("J132").Value = 5048.5884
$ws.Range("K132").Value = 5437.8138
$ws.Range("L132").Value = 15145.7652
$ws.Range("M132").Value = -2907.8138
$ws.Range("N132").Value = -20205.7652

# Row 136
$ws.Range("H136").Value = 2621.5818
$ws.Range("I136").Value = 2106.9412
$ws.Range("J136").Value = 3454.8096
$ws.Range("K136").Value = 6320.823600000001
$ws.Range("L136").Value = 10364.4288
$ws.Range("M136").Value = -3770.823600000001
$ws.Range("N136").Value = -15464.4288

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 1255.3529
$ws.Range("I80").Value = 2437.125
$ws.Range("J80").Value = 204.88889
$ws.Range("K80").Value = 2437.125
$ws.Range("L80").Value = 204.88889
$ws.Range("M80").Value = -1439.125
$ws.Range("N80").Value = -2200.88889

# Row 83
$ws.Range("H83").Value = 1255.3529
$ws.Range("I83").Value = 2437.125
$ws.Range("J83").Value = 204.88889
$ws.Range("K83").Value = 12185.625
$ws.Range("L83").Value = 1024.44445
$ws.Range("M83").Value = -7193.625
$ws.Range("N83").Value = -11008.44445

$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 3433.3333
$ws.Range("I10").Value = 2650
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 2650
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = -2511
$ws.Range("N10").Value = -5278

# Row 107
$ws.Range("H107").Value = 2016598.1
$ws.Range("I107").Value = 3676698.8
$ws.Range("J107").Value = 761.6429000000001
$ws.Range("K107").Value = 3676698.8
$ws.Range("L107").Value = 761.6429000000001
$ws.Range("M107").Value = -3674778.8
$ws.Range("N107").Value = -4601.6429

# Row 132
$ws.Range("H132").Value = 25363618
$ws.Range("I132").Value = 28572526
$ws.Range("J132").Value = 15153451
$ws.Range("K132").Value = 85717578
$ws.Range("L132").Value = 45460353
$ws.Range("M132").Value = -85715048
$ws.Range("N132").Value = -45465413

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1107.7222
$ws.Range("J113").Value = 1347.4166
$ws.Range("L113").Value = 4042.2498
$ws.Range("N113").Value = -8382.2498

# Row 122
$ws.Range("H122").Value = 2866.5908
$ws.Range("J122").Value = 6182.4736
$ws.Range("L122").Value = 55642.2624
$ws.Range("N122").Value = -60542.2624

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8282.5
$ws.Range("I70").Value = 9028.956
$ws.Range("J70").Value = 5829.857
$ws.Range("K70").Value = 9028.956
$ws.Range("L70").Value = 5829.857
$ws.Range("M70").Value = -8758.956
$ws.Range("N70").Value = -6369.857

# Row 73
$ws.Range("H73").Value = 8282.5
$ws.Range("I73").Value = 9028.956
$ws.Range("J73").Value = 5829.857
$ws.Range("K73").Value = 9028.956
$ws.Range("L73").Value = 5829.857
$ws.Range("M73").Value = -8092.956
$ws.Range("N73").Value = -7701.857

# Row 113
$ws.Range("H113").Value = 2480.3684
$ws.Range("I113").Value = 2537.3845
$ws.Range("K113").Value = 2537.3845
$ws.Range("M113").Value = -367.3845000000001

# Row 121
$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -53494

# Row 123
$ws.Range("H123").Value = 10325.7
$ws.Range("J123").Value = 10325.7
$ws.Range("L123").Value = 10325.7
$ws.Range("N123").Value = -15225.7

$ws = $wb.Worksheets.Item("LTW")
# Row 135
$ws.Range("H135").Value = 59239.5
$ws.Range("J135").Value = 59239.5
$ws.Range("L135").Value = 59239.5
$ws.Range("N135").Value = -69379.5

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 2220.5
$ws.Range("I136").Value = 1723.8889
$ws.Range("K136").Value = 5171.6667
$ws.Range("M136").Value = -2621.6667
